# Applies the "cryptos list" data refresh described in the commit
# "Updated cryptos list on Thu Aug  1 21:27:47 UTC 2024 with GitHub Actions".
# Column D (Price) and column E (Volume(1h)) are plain text cells (not numbers),
# so values that look numeric must be forced to stay text -- otherwise the COM
# layer auto-converts a bare "5.29" style string into a real number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Writes $Value into $CellRef while preserving its original text (General)
# formatting: briefly force a text NumberFormat so Excel does not reinterpret
# a numeric-looking string as a number/date, then restore the default style so
# no stray cell formatting is left behind.
function Set-TextValue {
    param($CellRef, $Value)
    $c = $ws.Range($CellRef)
    $c.NumberFormat = "@"
    $c.Value = $Value
    $c.Style = "Normal"
}

# Row 2 (Bitcoin)
$ws.Range("D2").Value = "64.853.21"
$ws.Range("E2").Value = "  -0.22%  "

# Row 3 (Ethereum)
$ws.Range("D3").Value = "3.165.60"
$ws.Range("E3").Value = "  -2.30%  "

# Row 4 (TetherUSD)
$ws.Range("E4").Value = "  -0.05%  "

# Row 5 (BNB)
Set-TextValue "D5" "571.71"
$ws.Range("E5").Value = "  -1.21%  "

# Row 6 (Solana)
Set-TextValue "D6" "165.60"
$ws.Range("E6").Value = "  -5.07%  "

# Row 7 (XRP)
Set-TextValue "D7" "0.591"
$ws.Range("E7").Value = "  -6.25%  "

# Row 8 (USDC)
$ws.Range("E8").Value = "  +0.04%  "

# Row 9 (Dogecoin)
$ws.Range("E9").Value = "  -3.90%  "

# Row 10 (Toncoin)
Set-TextValue "D10" "6.68"
$ws.Range("E10").Value = "  -1.50%  "

# Row 11 (Cardano)
Set-TextValue "D11" "0.386"
$ws.Range("E11").Value = "  -1.48%  "

# Row 12 (WrappedliquidstakedEther2.0)
$ws.Range("D12").Value = "3.713.11"
$ws.Range("E12").Value = "  -2.38%  "

# Row 13 (TRON)
$ws.Range("E13").Value = "  -0.71%  "

# Row 14 (WrappedBTC)
$ws.Range("D14").Value = "64.718.53"
$ws.Range("E14").Value = "  -0.67%  "

# Row 15 (Avalanche)
Set-TextValue "D15" "25.39"
$ws.Range("E15").Value = "  -1.03%  "

# Row 16 (WrappedEther)
$ws.Range("D16").Value = "3.156.73"
$ws.Range("E16").Value = "  -2.53%  "

# Row 17 (ShibaInu)
$ws.Range("E17").Value = "  -2.20%  "

# Row 18 (BitcoinCash)
Set-TextValue "D18" "413.83"
$ws.Range("E18").Value = "  -0.12%  "

# Row 19: Chainlink/Polkadot swap places in the ranking -> row 19 becomes Polkadot
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D19" "5.29"
$ws.Range("E19").Value = "  -2.03%  "

# Row 20: ... and row 20 becomes Chainlink
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D20" "12.63"
$ws.Range("E20").Value = "  -2.04%  "

# Row 21 (Uniswap)
Set-TextValue "D21" "7.12"
$ws.Range("E21").Value = "  -1.29%  "

# Row 22 (Dai)
$ws.Range("E22").Value = "  +0.11%  "

# Row 23 (Litecoin)
Set-TextValue "D23" "68.26"
$ws.Range("E23").Value = "  -3.03%  "

# Row 24 (Kaspa)
$ws.Range("E24").Value = "  -1.69%  "

# Row 25 (Polygon)
Set-TextValue "D25" "0.486"
$ws.Range("E25").Value = "  -2.01%  "

# Row 26 (PEPE)
$ws.Range("E26").Value = "  -6.66%  "

# Row 27 (InternetComputer(DFINITY))
Set-TextValue "D27" "8.89"
$ws.Range("E27").Value = "  -2.63%  "

# Row 28 (Binance-PegBSC-USD)
Set-TextValue "D28" "0.999"
$ws.Range("E28").Value = "  +0.00%  "

# Row 29 (PancakeSwap)
$ws.Range("E29").Value = "  -3.12%  "

# Row 30 (EthereumClassic)
Set-TextValue "D30" "21.30"
$ws.Range("E30").Value = "  -2.49%  "

# Row 31 (NEARProtocol)
Set-TextValue "D31" "4.91"
$ws.Range("E31").Value = "  -1.97%  "

# Row 32 (Aptos)
$ws.Range("E32").Value = "  -1.66%  "

# Row 33 (Fetch.AI)
$ws.Range("E33").Value = "  -2.62%  "

# Row 34 (Monero)
Set-TextValue "D34" "155.80"
$ws.Range("E34").Value = "  -0.69%  "

# Row 35 (ImmutableX)
$ws.Range("E35").Value = "  -3.33%  "

# Row 36 (Maker)
$ws.Range("D36").Value = "2.716.26"
$ws.Range("E36").Value = "  -4.17%  "

# Row 37 (Stacks)
$ws.Range("E37").Value = "  -3.47%  "

# Row 38 (EnergySwap)
Set-TextValue "D38" "23.60"
$ws.Range("E38").Value = "  -7.75%  "

# Row 39 (Filecoin)
$ws.Range("E39").Value = "  -3.50%  "

# Row 40 (Mantle)
$ws.Range("E40").Value = "  -4.25%  "

# Row 41 (Hedera)
Set-TextValue "D41" "0.0631"
$ws.Range("E41").Value = "  +0.27%  "

# Row 42 (RenderToken)
$ws.Range("E42").Value = "  -4.42%  "

# Row 43 (VeChain)
$ws.Range("E43").Value = "  -0.94%  "

# Row 44 (Bittensor)
Set-TextValue "D44" "288.90"
$ws.Range("E44").Value = "  -5.97%  "

# Row 45 (InjectiveProtocol)
Set-TextValue "D45" "21.15"
$ws.Range("E45").Value = "  -4.73%  "

# Row 46 (FirstDigitalUSD)
Set-TextValue "D46" "0.999"
$ws.Range("E46").Value = "  -0.10%  "

# Row 47 (Stellar)
$ws.Range("E47").Value = "  -3.14%  "

# Row 48 (dogwifhat)
$ws.Range("E48").Value = "  -11.29%  "

# Row 49 (WhiteBITCoin)
Set-TextValue "D49" "10.47"
$ws.Range("E49").Value = "  +0.74%  "

# Row 50 (Cosmos)
Set-TextValue "D50" "5.75"
$ws.Range("E50").Value = "  -1.18%  "

# Row 51 (ONDO)
Set-TextValue "D51" "0.894"
$ws.Range("E51").Value = "  -4.90%  "
